$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipment")

$data = @{
    10 = @("PT-8", "PDT-6", 3, 4, 4, 4, 8)
    11 = @("PT-7", "PDT-6", 4, 4, 2, 4, 4)
    12 = @("PT-6", "PDT-4", 9, 7, 4, 2, 4)
    13 = @("PT-7", "PDT-2", 4, 2, 6, 6, 8)
    14 = @("PT-6", "PDT-5", 6, 9, 9, 5, 6)
    15 = @("PT-4", "PDT-6", 5, 6, 3, 5, 4)
    16 = @("PT-2", "PDT-7", 7, 2, 7, 5, 4)
    17 = @("PT-5", "PDT-6", 5, 3, 5, 7, 2)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
    $ws.Cells.Item($row, 8).Value = $values[6]
}
